$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# 9x39 AP and DMG boost: bump the "ammo_k_ap" and "ammo_k_hit" values for
# ammo_9x39_pab9 (row 21) and ammo_9x39_ap (row 22)
$ws.Range("G21").Value = 0.34
$ws.Range("H21").Value = 1.04

$ws.Range("G22").Value = 0.55
$ws.Range("H22").Value = 1.04

$excel.Calculate()

# Update the active selection to match the edited area
$ws.Activate()
$ws.Range("H20").Select()
